# "Final bulk import: > separator, Excel input, duplicate detection, mixed categories"
#
# - Category labels use a plain "Health > Sleep" / "Health > Daily Wellness"
#   separator instead of the old "Health -> Sleep" style arrow glyph.
# - Row 3 date corrected (now a duplicate-detection scenario: different
#   dates for entries in the same category) and rows 4/5 dates likewise
#   adjusted.
# - The two example rows exercising "Training > Cardio > Upper Body" /
#   Excel-date-input placeholders (rows 6 & 7) are removed — the template
#   now ships with the 4 "real" example rows only.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the category column to the new "Health > ..." text (replaces the
# old unicode-arrow separator) for the two Health categories still present.
$ws.Range("A2").Value = "Health > Sleep"
$ws.Range("A3").Value = "Health > Sleep"
$ws.Range("A4").Value = "Health > Daily Wellness"
$ws.Range("A5").Value = "Health > Daily Wellness"

# Correct the sample dates (rows 3-5) to distinct days.
$ws.Range("B3").Value = 45971
$ws.Range("B4").Value = 45972
$ws.Range("B5").Value = 45970

# Drop the two "Training > Cardio > Upper Body" example rows entirely.
$ws.Rows("6:7").Delete()

# Match the author's final selection state.
$ws.Range("A7").Select()
